# Update "想去人数" (want-to-go count) values that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览": F3 269 -> 270, F4 922 -> 923
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 270
$ws1.Range("F4").Value = 923

# Sheet "全部类型": F4 269 -> 270, F5 922 -> 923
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 270
$ws4.Range("F5").Value = 923
